$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (col A and col B)
$ws.Columns.Item(1).ColumnWidth = 14.75
$ws.Columns.Item(2).ColumnWidth = 14.584

# Update cell values
$ws.Cells.Item(1, 1).Value = -0.36992273206928417
$ws.Cells.Item(1, 2).Value = 0.36888264975718243
$ws.Cells.Item(2, 1).Value = -0.281619912327832
$ws.Cells.Item(2, 2).Value = 0.2783293756959466
$ws.Cells.Item(3, 1).Value = -0.175376036625849
$ws.Cells.Item(3, 2).Value = 0.1743540213267316
$ws.Cells.Item(4, 1).Value = -0.16235402144005207
$ws.Cells.Item(4, 2).Value = 0.16143260611259613
$ws.Cells.Item(5, 1).Value = -0.155432606526694
$ws.Cells.Item(5, 2).Value = 0.15358399000643175
$ws.Cells.Item(6, 1).Value = -0.08285086169964995
$ws.Cells.Item(6, 2).Value = 0.08275270041032412
$ws.Cells.Item(7, 1).Value = -0.062752700916473
$ws.Cells.Item(7, 2).Value = 0.06254455948053739
$ws.Cells.Item(8, 1).Value = -0.042544559990777664
$ws.Cells.Item(8, 2).Value = 0.0424107020348119
$ws.Cells.Item(9, 1).Value = -0.03641070247176792
$ws.Cells.Item(9, 2).Value = 0.03631145748741549
$ws.Cells.Item(10, 1).Value = -0.030311457927894025
$ws.Cells.Item(10, 2).Value = 0.030302257998577886
$ws.Cells.Item(11, 1).Value = -0.0258022584304598
$ws.Cells.Item(11, 2).Value = 0.02578704452907843
$ws.Cells.Item(12, 1).Value = -0.019787044970023704
$ws.Cells.Item(12, 2).Value = 0.019747303228684565
$ws.Cells.Item(13, 1).Value = -0.013747303672284161
$ws.Cells.Item(13, 2).Value = 0.013739890753538653
$ws.Cells.Item(14, 1).Value = -0.0017398912309367631
$ws.Cells.Item(14, 2).Value = 0.0017343817412802665
$ws.Cells.Item(15, 1).Value = -0.021053892452812306
$ws.Cells.Item(15, 2).Value = 0.02102809849884224
$ws.Cells.Item(16, 1).Value = -0.015028098944527724
$ws.Cells.Item(16, 2).Value = 0.015004452756143927
$ws.Cells.Item(17, 1).Value = -0.009004453203687035
$ws.Cells.Item(17, 2).Value = 0.008999999535381775
$ws.Cells.Item(18, 1).Value = -0.07726273784150806
$ws.Cells.Item(18, 2).Value = 0.07718605347760743
$ws.Cells.Item(19, 1).Value = -0.068186053885265
$ws.Cells.Item(19, 2).Value = 0.0676086112933807
$ws.Cells.Item(20, 1).Value = -0.018013865171068844
$ws.Cells.Item(20, 2).Value = 0.018004324470101807
$ws.Cells.Item(21, 1).Value = -0.009004324888271853
$ws.Cells.Item(21, 2).Value = 0.008999999581416063
$ws.Cells.Item(22, 1).Value = -0.09395333948996765
$ws.Cells.Item(22, 2).Value = 0.09363841827550701
$ws.Cells.Item(23, 1).Value = -0.08463841869911182
$ws.Cells.Item(23, 2).Value = 0.08412761109716715
$ws.Cells.Item(24, 1).Value = -0.042127611709958934
$ws.Cells.Item(24, 2).Value = 0.041999999383934394
$ws.Cells.Item(25, 1).Value = -0.0647331287339874
$ws.Cells.Item(25, 2).Value = 0.06461796864612523
$ws.Cells.Item(26, 1).Value = -0.05861796907429451
$ws.Cells.Item(26, 2).Value = 0.05847453082607146
$ws.Cells.Item(27, 1).Value = -0.038414033178322704
$ws.Cells.Item(27, 2).Value = 0.038076058700554416
$ws.Cells.Item(28, 1).Value = -0.0320760591342415
$ws.Cells.Item(28, 2).Value = 0.03185513618953184
$ws.Cells.Item(29, 1).Value = -0.019855136658938122
$ws.Cells.Item(29, 2).Value = 0.01976362687358524
$ws.Cells.Item(30, 1).Value = 0.00023637261192899217
$ws.Cells.Item(30, 2).Value = -0.00024704742201953067
$ws.Cells.Item(31, 1).Value = 0.015247046935289887
$ws.Cells.Item(31, 2).Value = -0.015259266063043952
$ws.Cells.Item(32, 1).Value = -0.021387101311943546
$ws.Cells.Item(32, 2).Value = 0.021368926733398297
